$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window view tweak (bookViews/workbookView) ---
$excel.Width = 1000
$excel.Height = 700

# --- New testplan section: "Entity report patient" ---
# Row 81: section header (style copied from an existing header row, e.g. row 75)
$ws.Range("A75:B75").Copy() | Out-Null
$ws.Range("A81:B81").PasteSpecial(-4122) | Out-Null
$ws.Range("A81").Value = "Entity report patient"

# Row 82
$ws.Range("A79").Copy() | Out-Null
$ws.Range("A82").PasteSpecial(-4122) | Out-Null
$ws.Range("A82").Value = "Go to the data tab"

# Row 83
$ws.Range("A79").Copy() | Out-Null
$ws.Range("A83").PasteSpecial(-4122) | Out-Null
$ws.Range("A83").Value = "Click on the entityreport button "

# Row 84
$ws.Range("A80:B80").Copy() | Out-Null
$ws.Range("A84:B84").PasteSpecial(-4122) | Out-Null
$ws.Range("A84").Value = "Try to search in the table"
$ws.Range("B84").Value = "A patient view like screen should show"

# Row 85
$ws.Range("A80:B80").Copy() | Out-Null
$ws.Range("A85:B85").PasteSpecial(-4122) | Out-Null
$ws.Range("A85").Value = "Try to download "
$ws.Range("B85").Value = "Does it work?"

# Row 86
$ws.Range("A79").Copy() | Out-Null
$ws.Range("A86").PasteSpecial(-4122) | Out-Null
$ws.Range("A86").Value = "Does a file download with the name: patient_data_*patient name*.pdf?"

# Row 87
$ws.Range("A80:B80").Copy() | Out-Null
$ws.Range("A87:B87").PasteSpecial(-4122) | Out-Null
$ws.Range("A87").Value = "Check the file"
$ws.Range("B87").Value = "Does it contain the content of the patient table? (yes genotype shows very ugly)"

# Row 88
$ws.Range("A80:B80").Copy() | Out-Null
$ws.Range("A88:B88").PasteSpecial(-4122) | Out-Null
$ws.Range("A88").Value = "Now once go the the patient view again and select the same patient as you selected in the entity report"
$ws.Range("B88").Value = "Does everything appear the same as before?"

$excel.CutCopyMode = $false

# --- Update the view to reflect the new bottom of the sheet ---
$ws.Range("A72").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 72
$ws.Range("B97").Select() | Out-Null
